# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and
# the blank/page-break paragraphs that immediately surround it, right
# after the "LOB1036: Geometria Analítica (Requisito fraco)" paragraph.

$d = $word.ActiveDocument

$anchorText = "LOB1036: Geometria Analítica (Requisito fraco)"
$targetText = "Ver no Jupiter Salvar em pdf Salvar em docx"

$anchorIndex = -1
$targetIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains($anchorText)) {
        $anchorIndex = $i
    }
    if ($t.Contains($targetText)) {
        $targetIndex = $i
    }
}

if ($anchorIndex -eq -1 -or $targetIndex -eq -1) {
    throw "Could not locate anchor/target paragraphs (anchor=$anchorIndex target=$targetIndex)"
}

# The four paragraphs to delete are: the blank paragraph right after the
# anchor, the "Ver no Jupiter..." paragraph, the blank paragraph after
# it, and the following page-break paragraph (with jc=left) that
# precedes the trailing blank/page-break paragraphs kept in the doc.
$firstToRemove = $anchorIndex + 1
$lastToRemove = $targetIndex + 2

$startRange = $d.Paragraphs.Item($firstToRemove).Range.Start
$endRange = $d.Paragraphs.Item($lastToRemove).Range.End

$r = $d.Range($startRange, $endRange)
$r.Delete()
